# Semana 04 de 2025: refresh the poisson-by-event table with the new
# epidemiological week data (new/changed events, recomputed Esperado/
# Observado/valor p, and several retired events removed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previous week's data rows (2-35); header row 1 is untouched.
$ws.Range("A2:E35").Clear()

# evento, nom_eve, Esperado, Observado, "valor p"
$data = @(
    ,@("113", "Desnutrici”n aguda en menores de 5 anos", 2, 1, 0.27)
    ,@("115", "Cancer en menores de 18 anos", 1, 1, 0.37)
    ,@("155", "Cancer de la mama y cuello uterino", 7, 2, 0.02)
    ,@("210", "Dengue", 1, 28, 0)
    ,@("215", "Defectos congenitos", 1, 2, 0.18)
    ,@("220", "Dengue grave", 0, 0, 1)
    ,@("300", "Agresiones por animales potencialmente transmisores de rabia", 50, 40, 0.02)
    ,@("340", "Hepatitis b, c y coinfeccion hepatitis b y delta", 1, 0, 0.37)
    ,@("342", "Enfermedades huerfanas - raras", 2, 3, 0.18)
    ,@("346", "Ira por virus nuevo", 294, 1, 0)
    ,@("348", "Infeccion respiratoria aguda grave irag inusitada", 1, 0, 0.37)
    ,@("355", "Enfermedad transmitida por alimentos o agua (eta)", 0, 0, 1)
    ,@("356", "Intento de suicidio", 11, 6, 0.04)
    ,@("357", "Iad - infecciones asociadas a dispositivos - individual", 1, 1, 0.37)
    ,@("365", "Intoxicaciones", 6, 1, 0.01)
    ,@("455", "Leptospirosis", 1, 1, 0.37)
    ,@("465", "Malaria", 0, 3, 0)
    ,@("549", "Morbilidad materna extrema", 6, 3, 0.09)
    ,@("560", "Mortalidad perinatal y neonatal tardia", 1, 0, 0.37)
    ,@("620", "Parotiditis", 1, 1, 0.37)
    ,@("740", "Sifilis congenita", 0, 0, 1)
    ,@("750", "Sifilis gestacional", 2, 2, 0.27)
    ,@("813", "Tuberculosis", 8, 9, 0.12)
    ,@("831", "Varicela individual", 5, 4, 0.18)
    ,@("850", "Vih/sida/mortalidad por sida", 7, 2, 0.02)
)

$firstRow = 2
$lastRow = $firstRow + $data.Count - 1

# Column A codes (e.g. "113", "300") look numeric but must stay TEXT,
# matching the source data. Pre-formatting the whole column range as
# Text keeps Excel from auto-converting them to numbers; ClearFormats
# afterwards drops that temporary number-format again so cells end up
# with no explicit style, same as the original rows.
$colA = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, 1))
$colA.NumberFormat = "@"

$r = $firstRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$colA.ClearFormats()

Write-Output "Updated $($data.Count) rows (${firstRow} to ${lastRow})"
